$d = $word.ActiveDocument

# 1. Fix "Bhandri" -> "Bhandari"
$d.Content.Find.Execute("Bhandri", $false, $false, $false, $false, $false, $true, 1, $false, "Bhandari", 2)

# 2. "Android and iOS platforms" -> "Windows Laptop or PC" (Bound scope sentence)
$d.Content.Find.Execute("This guide applies to Android and iOS platforms for users engaging with the game.", $false, $false, $false, $false, $false, $true, 1, $false, "This guide applies to Windows Laptop or PC for users engaging with the game.", 2)

# 3. Quick start guide sentence - merge runs (no textual change needed, but ensure text consistent)
$d.Content.Find.Execute("Quick start guide to begin playing in 3 steps", $false, $false, $false, $false, $false, $true, 1, $false, "Quick start guide to begin playing in 3 steps", 2)

# 4. Trophy / Score / Replay / Menu - merge runs (no textual change needed)
$d.Content.Find.Execute("Trophy " + [char]8226 + " Score " + [char]8226 + " Replay  " + [char]8226 + " Menu", $false, $false, $false, $false, $false, $true, 1, $false, "Trophy " + [char]8226 + " Score " + [char]8226 + " Replay  " + [char]8226 + " Menu", 2)

# 5. "Tested on Android 11 tablets and iPadOS 17+" -> "Tested on Windows Laptop or PC"
$d.Content.Find.Execute("Tested on Android 11 tablets and iPadOS 17+", $false, $false, $false, $false, $false, $true, 1, $false, "Tested on Windows Laptop or PC", 2)
